$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New species-code values being filled into column B (SPECIES_CODE) for rows
# that previously only had column A (OLD_MVS_CODE) populated.
$ws.Range("B54").Value  = "LARV"
$ws.Range("B92").Value  = "SPH_BORE"
$ws.Range("B104").Value = "SYN_FLOR"
$ws.Range("B105").Value = "SYN_LOUI"
$ws.Range("B107").Value = "TRI_INSC"

# Reposition the view/selection to match where the author ended up working.
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B110").Select()
